# Update "想去人数" (F) and "最低票价" (G) values on the 展览, 演出 and 全部类型 sheets
# per the output-data refresh commit (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 3850
$ws.Range("F4").Value = 1377
$ws.Range("F5").Value = 3843
$ws.Range("G6").Value = 65
$ws.Range("F7").Value = 213
$ws.Range("F8").Value = 62
$ws.Range("F9").Value = 8837
$ws.Range("F13").Value = 353
$ws.Range("F17").Value = 387
$ws.Range("F18").Value = 11232
$ws.Range("F21").Value = 78
$ws.Range("F25").Value = 145
$ws.Range("F33").Value = 2085
$ws.Range("F37").Value = 920
$ws.Range("F38").Value = 4104
$ws.Range("F39").Value = 2573
$ws.Range("F41").Value = 2590
$ws.Range("F45").Value = 74
$ws.Range("F46").Value = 362
$ws.Range("F47").Value = 52

# --- Sheet: 演出 (Performance) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("G6").Value = 108
$ws.Range("F9").Value = 52
$ws.Range("F10").Value = 34
$ws.Range("F18").Value = 179
$ws.Range("F20").Value = 32

# --- Sheet: 全部类型 (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 3850
$ws.Range("F5").Value = 3843
$ws.Range("G6").Value = 65
$ws.Range("G7").Value = 108
$ws.Range("F8").Value = 52
$ws.Range("F9").Value = 213
$ws.Range("F10").Value = 62
$ws.Range("F11").Value = 8837
$ws.Range("F12").Value = 34
$ws.Range("F16").Value = 353
$ws.Range("F20").Value = 387
$ws.Range("F21").Value = 11232
$ws.Range("F24").Value = 145
$ws.Range("F30").Value = 2085
$ws.Range("F34").Value = 920
$ws.Range("F35").Value = 179
$ws.Range("F37").Value = 2573
$ws.Range("F39").Value = 2590
$ws.Range("F44").Value = 74
$ws.Range("F45").Value = 362
$ws.Range("F47").Value = 52
